$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.837082743644714
$ws.Range("B1").Value = -1
$ws.Range("C1").Value = 2.451897859573364
$ws.Range("D1").Value = 1.239700555801392
$ws.Range("E1").Value = 0.8814833164215088
